# Add gems reward desc - populate columns M (cn) / N (en) on the "宝石"
# (gems) sheet with first-time-purchase reward descriptions, plus the
# cosmetic view tweaks captured in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row
$ws.Range("M1").Value = "cn"
$ws.Range("N1").Value = "en"

# Reward descriptions per gem pack (rows 2-6), Chinese + English pairs
$ws.Range("M2").Value = "首次购买赠送10万金币"
$ws.Range("N2").Value = "First-time-purchase: 100k gold"

$ws.Range("M3").Value = "首次购买赠送白色霸王龙*1"
$ws.Range("N3").Value = "First-time-purchase: white Tyrent x 1"

$ws.Range("M4").Value = "首次购买赠送绿色霸王龙*1"
$ws.Range("N4").Value = "First-time-purchase: green Tyrent x 1"

$ws.Range("M5").Value = "首次购买赠送蓝色霸王龙*1和蓝色震龙*1"
$ws.Range("N5").Value = "First-time-purchase: blue Tyrent x 1, blue earthquake x 1"

$ws.Range("M6").Value = "首次购买赠送紫色霸王龙*1"
$ws.Range("N6").Value = "First-time-purchase: purple Tyrent x 2"

# Match styling: M column centered horizontally+vertically, N column
# vertically centered (reuses the workbook's existing alignment styles).
$ws.Range("M2:M6").HorizontalAlignment = -4108
$ws.Range("M2:M6").VerticalAlignment = -4108
$ws.Range("N2:N6").VerticalAlignment = -4108

# Column widths for the two new columns
$ws.Columns.Item(13).ColumnWidth = 29.5
$ws.Columns.Item(14).ColumnWidth = 57.5

# Selection / view state tweaks from the commit
$ws.Select()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H18").Select()
